# Remove the 2025-09-12 placeholder row (row 2) from the "Chart" sheet.
# That row had no "Not indexed" / "Indexed" data yet (blank shared strings),
# so the whole row is deleted and everything below shifts up by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")
$ws.Rows.Item(2).Delete()
